$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 420.32144
$ws.Range("I17").Value = 1500
$ws.Range("J17").Value = 240.375
$ws.Range("K17").Value = 4500
$ws.Range("L17").Value = 721.125
$ws.Range("M17").Value = -4332
$ws.Range("N17").Value = -1057.125

$ws.Range("H18").Value = 600
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 600
$ws.Range("M18").Value = ""
$ws.Range("N18").Value = -1168

$ws.Range("H34").Value = 3701.4546
$ws.Range("I34").Value = 2421.6
$ws.Range("J34").Value = 16500
$ws.Range("K34").Value = 2421.6
$ws.Range("L34").Value = 16500
$ws.Range("M34").Value = -2218.6
$ws.Range("N34").Value = -16906

$ws.Range("H36").Value = 3701.4546
$ws.Range("I36").Value = 2421.6
$ws.Range("J36").Value = 16500
$ws.Range("K36").Value = 2421.6
$ws.Range("L36").Value = 16500
$ws.Range("M36").Value = -1706.6
$ws.Range("N36").Value = -17930

$ws.Range("H120").Value = 50000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 50000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 50000
$ws.Range("N120").Value = -59676

$ws.Range("H123").Value = 40998.75
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 40998.75
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 40998.75
$ws.Range("N123").Value = -50798.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 50000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 50000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 50000
$ws.Range("N113").Value = -58678

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = ""

$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = ""

$ws.Range("H117").Value = 33333
$ws.Range("I117").Value = 33333
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 33333
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = -28744

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 708.08
$ws.Range("I22").Value = 447.17648
$ws.Range("J22").Value = 1262.5
$ws.Range("K22").Value = 447.17648
$ws.Range("L22").Value = 1262.5
$ws.Range("M22").Value = -97.17648000000003
$ws.Range("N22").Value = -1962.5

$ws.Range("H31").Value = 5025.156
$ws.Range("I31").Value = 1622.4117
$ws.Range("J31").Value = 15542.728
$ws.Range("K31").Value = 1622.4117
$ws.Range("L31").Value = 15542.728
$ws.Range("M31").Value = -1327.4117
$ws.Range("N31").Value = -16132.728

$ws.Range("H34").Value = 5025.156
$ws.Range("I34").Value = 1622.4117
$ws.Range("J34").Value = 15542.728
$ws.Range("K34").Value = 1622.4117
$ws.Range("L34").Value = 15542.728
$ws.Range("M34").Value = -1420.4117
$ws.Range("N34").Value = -15946.728

$ws.Range("H88").Value = 41000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 41000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 41000
$ws.Range("N88").Value = -41812

$ws.Range("H91").Value = 41000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 41000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 41000
$ws.Range("N91").Value = -43808

$ws.Range("H97").Value = 32423.75
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 32423.75
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 32423.75
$ws.Range("N97").Value = -34405.75

$ws.Range("H98").Value = 47249
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 47249
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 47249
$ws.Range("N98").Value = -51741

$ws.Range("H104").Value = 37039.668
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 37039.668
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 37039.668
$ws.Range("N104").Value = -42281.668

$ws.Range("H105").Value = 1573.762
$ws.Range("I105").Value = 1449.9445
$ws.Range("J105").Value = 2316.6667
$ws.Range("K105").Value = 1449.9445
$ws.Range("L105").Value = 2316.6667
$ws.Range("M105").Value = 297.0554999999999
$ws.Range("N105").Value = -5810.6667

$ws.Range("H106").Value = 35167.75
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 35167.75
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 35167.75
$ws.Range("M106").Value = ""
$ws.Range("N106").Value = -37691.75

$ws.Range("H107").Value = 1039
$ws.Range("I107").Value = 734.6
$ws.Range("J107").Value = 1800
$ws.Range("K107").Value = 734.6
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 1185.4
$ws.Range("N107").Value = -5640

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""

$ws.Range("H134").Value = 5104.3438
$ws.Range("I134").Value = 5570.7915
$ws.Range("J134").Value = 3705
$ws.Range("K134").Value = 16712.3745
$ws.Range("L134").Value = 11115
$ws.Range("M134").Value = -14177.3745
$ws.Range("N134").Value = -16185

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1485543
$ws.Range("I5").Value = 1688.5
$ws.Range("J5").Value = 2672626.5
$ws.Range("K5").Value = 5065.5
$ws.Range("L5").Value = 8017879.5
$ws.Range("M5").Value = -4953.5
$ws.Range("N5").Value = -8018103.5

$ws.Range("H122").Value = 2894.6956
$ws.Range("I122").Value = 1108.5
$ws.Range("J122").Value = 3270.7368
$ws.Range("K122").Value = 9976.5
$ws.Range("L122").Value = 29436.6312
$ws.Range("M122").Value = -7526.5
$ws.Range("N122").Value = -34336.6312

$ws.Range("H132").Value = 2283.423
$ws.Range("I132").Value = 934.4545000000001
$ws.Range("J132").Value = 3272.6667
$ws.Range("K132").Value = 8410.0905
$ws.Range("L132").Value = 29454.0003
$ws.Range("M132").Value = -5880.0905
$ws.Range("N132").Value = -34514.0003

$ws.Range("H135").Value = 1485543
$ws.Range("I135").Value = 1688.5
$ws.Range("J135").Value = 2672626.5
$ws.Range("K135").Value = 15196.5
$ws.Range("L135").Value = 24053638.5
$ws.Range("M135").Value = -12661.5
$ws.Range("N135").Value = -24058708.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 57857.5
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 57857.5
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 57857.5
$ws.Range("N36").Value = -58981.5

$ws.Range("H122").Value = 6949
$ws.Range("I122").Value = 5937.6665
$ws.Range("J122").Value = 11500
$ws.Range("K122").Value = 17812.9995
$ws.Range("L122").Value = 34500
$ws.Range("M122").Value = -15362.9995
$ws.Range("N122").Value = -39400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11637.692
$ws.Range("I62").Value = 3823.75
$ws.Range("J62").Value = 24140
$ws.Range("K62").Value = 3823.75
$ws.Range("L62").Value = 24140
$ws.Range("M62").Value = -3199.75
$ws.Range("N62").Value = -25388

$ws.Range("H65").Value = 11637.692
$ws.Range("I65").Value = 3823.75
$ws.Range("J65").Value = 24140
$ws.Range("K65").Value = 19118.75
$ws.Range("L65").Value = 120700
$ws.Range("M65").Value = -15998.75
$ws.Range("N65").Value = -126940

$ws.Range("H80").Value = 39800
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 39800
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 39800
$ws.Range("N80").Value = -41796

$ws.Range("H81").Value = 2655.5
$ws.Range("I81").Value = 2624.182
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 5248.364
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -4187.364
$ws.Range("N81").Value = -8122

$ws.Range("H82").Value = 32550
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 32550
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 32550
$ws.Range("N82").Value = -33316

$ws.Range("H83").Value = 39800
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 39800
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 119400
$ws.Range("N83").Value = -129384

$ws.Range("H84").Value = 2655.5
$ws.Range("I84").Value = 2624.182
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 26241.82
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -20937.82
$ws.Range("N84").Value = -40608

$ws.Range("H85").Value = 32550
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 32550
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 32550
$ws.Range("N85").Value = -35202

$ws.Range("H107").Value = 614.2857
$ws.Range("I107").Value = 331.33334
$ws.Range("J107").Value = 826.5
$ws.Range("K107").Value = 994.0000200000001
$ws.Range("L107").Value = 2479.5
$ws.Range("M107").Value = 925.9999799999999
$ws.Range("N107").Value = -6319.5
